# Sales report PDF format bug fix:
# A new order row was missing from the top of the data table. Insert it
# as the new row 2 (right after the header row), pushing the existing
# data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new blank row at row 2 (above the former first data row),
# shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new order's data.
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = 6
$ws.Cells.Item(2, 3).Value = 52
# The Date column stores its values as plain text (matching the rest of
# the column), so format the cell as text before writing it - otherwise
# Excel auto-converts the recognizable "yyyy-mm-dd" string into a date.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "2023-11-24"
$ws.Cells.Item(2, 5).Value = "COD"
$ws.Cells.Item(2, 6).Value = 53
$ws.Cells.Item(2, 7).Value = 1389
$ws.Cells.Item(2, 8).Value = "delivered"
